# Quantum Computing Inc. workbook cleanup:
# Collapse the D:I (open/close/high/low/shares/ticker) columns for
# rows 2-33 down to QUBT-only values (the extra ticker rows that had
# been pulled in from other files are normalized back to QUBT data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @(open, close, high, low, shares_outstanding)
$rowValues = @{
    2 = @(0.6000000238418579, 0.6000000238418579, 0.6000000238418579, 0.4000000059604645, 223934400)
    3 = @(0.800000011920929, 6, 12.60000038146973, 0.6000000238418579, 223934400)
    4 = @(15.80000019073486, 5.599999904632568, 15.80000019073486, 3.599999904632568, 223934400)
    5 = @(2.400000095367432, 5.199999809265137, 8.489999771118164, 0.0149999996647238, 223934400)
    6 = @(12.5, 7, 12.5, 6, 223934400)
    7 = @(3, 5.349999904632568, 5.349999904632568, 2.849999904632568, 223934400)
    8 = @(4.050000190734863, 3.5, 6, 3.5, 223934400)
    9 = @(4.849999904632568, 5, 7.5, 3, 223934400)
    10 = @(3.5, 3.5, 4.110000133514404, 2.200000047683716, 223934400)
    11 = @(2.950000047683716, 2.25, 3.25, 1.850000023841858, 223934400)
    12 = @(2.055000066757202, 1.5, 2.299999952316284, 1.289999961853027, 223934400)
    13 = @(2.680000066757202, 2.25, 2.869999885559082, 1.950000047683716, 223934400)
    14 = @(3.724999904632568, 3.414999961853028, 4.300000190734863, 3.380000114440918, 223934400)
    15 = @(15.15499973297119, 11.1899995803833, 17, 9.199999809265137, 223934400)
    16 = @(7.010000228881836, 6.110000133514404, 8.970000267028809, 6.099999904632568, 223934400)
    17 = @(5.224999904632568, 9.340000152587891, 10.43000030517578, 4.989999771118164, 223934400)
    18 = @(6.25, 6.25, 6.889999866485596, 5.03000020980835, 223934400)
    19 = @(3.480000019073486, 2.690000057220459, 3.75, 2.049999952316284, 223934400)
    20 = @(2.599999904632568, 1.529999971389771, 2.75, 1.470000028610229, 223934400)
    21 = @(2.400000095367432, 3.069999933242798, 3.440000057220459, 2.269999980926514, 223934400)
    22 = @(2.5, 2.230000019073486, 2.670000076293945, 1.835999965667725, 223934400)
    23 = @(1.570000052452087, 1.759999990463257, 1.970000028610229, 1.5, 223934400)
    24 = @(1.360000014305115, 1.149999976158142, 1.580000042915344, 1.129999995231628, 223934400)
    25 = @(1.210000038146973, 1.419999957084656, 1.490000009536743, 1.149999976158142, 223934400)
    26 = @(1.080000042915344, 0.7269999980926514, 1.090000033378601, 0.6340000033378601, 223934400)
    27 = @(0.9290000200271606, 0.8109999895095825, 0.959999978542328, 0.75, 223934400)
    28 = @(1.029999971389771, 0.7680000066757202, 1.029999971389771, 0.6800000071525574, 223934400)
    29 = @(0.5130000114440918, 0.6819999814033508, 0.8199999928474426, 0.3549999892711639, 223934400)
    30 = @(0.6520000100135803, 1.120000004768372, 1.649999976158142, 0.6389999985694885, 223934400)
    31 = @(16.69000053405762, 10.38000011444092, 20.36000061035156, 6.409999847412109, 223934400)
    32 = @(8.199999809265137, 6.739999771118164, 8.399999618530273, 5.760000228881836, 223934400)
    33 = @(18.72999954223633, 14.84000015258789, 21.8799991607666, 14.625, 223934400)
}

foreach ($row in $rowValues.Keys) {
    $vals = $rowValues[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
    $ws.Range("F$row").Value = $vals[2]
    $ws.Range("G$row").Value = $vals[3]
    $ws.Range("H$row").Value = $vals[4]
    $ws.Range("I$row").Value = "QUBT"
}
